$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value = "Day"
$ws.Cells.Item(1,2).Value = "Time"
$ws.Cells.Item(1,3).Value = "Module Code"
$ws.Cells.Item(1,4).Value = "Module Title"
$ws.Cells.Item(1,5).Value = "Hours"
$ws.Cells.Item(1,6).Value = "Class Type"
$ws.Cells.Item(1,7).Value = "Lecturer"
$ws.Cells.Item(1,8).Value = "Room"
$ws.Cells.Item(1,9).Value = "Block"
$ws.Cells.Item(1,10).Value = "Group"
$ws.Cells.Item(1,11).Value = "Level"
$ws.Cells.Item(1,12).Value = "Course"

# Row 2
$ws.Cells.Item(2,1).Value = "SUN"
$ws.Cells.Item(2,2).Value = "7:00-9:00"
$ws.Cells.Item(2,3).Value = "5CS022"
$ws.Cells.Item(2,4).Value = "Human Computer Interaction"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = "Tutorial"
$ws.Cells.Item(2,7).Value = "Mr. Apurba Neupane"
$ws.Cells.Item(2,8).Value = "SR-03 Wolves"
$ws.Cells.Item(2,9).Value = "WLV"
$ws.Cells.Item(2,10).Value = "L5CG10"
$ws.Cells.Item(2,11).Value = 5
$ws.Cells.Item(2,12).Value = "BCS"

# Row 3
$ws.Cells.Item(3,1).Value = "SUN"
$ws.Cells.Item(3,2).Value = "9:30-12:00"
$ws.Cells.Item(3,3).Value = "5CS024"
$ws.Cells.Item(3,4).Value = "Collaborative Development"
$ws.Cells.Item(3,5).Value = 2.5
$ws.Cells.Item(3,6).Value = "Workshop"
$ws.Cells.Item(3,7).Value = "Mr. Anmol Adhikari"
$ws.Cells.Item(3,8).Value = "SR-03 Wolves"
$ws.Cells.Item(3,9).Value = "WLV"
$ws.Cells.Item(3,10).Value = "L5CG10"
$ws.Cells.Item(3,11).Value = 5
$ws.Cells.Item(3,12).Value = "BCS"

# Row 4
$ws.Cells.Item(4,1).Value = "MON"
$ws.Cells.Item(4,2).Value = "7:00-9:30"
$ws.Cells.Item(4,3).Value = "5CS022"
$ws.Cells.Item(4,4).Value = "Human Computer Interaction"
$ws.Cells.Item(4,5).Value = 2.5
$ws.Cells.Item(4,6).Value = "Workshop"
$ws.Cells.Item(4,7).Value = "Mr. Apurba Neupane"
$ws.Cells.Item(4,8).Value = "Lab-03 Gahanapokhari"
$ws.Cells.Item(4,9).Value = "HCK"
$ws.Cells.Item(4,10).Value = "L5CG10"
$ws.Cells.Item(4,11).Value = 5
$ws.Cells.Item(4,12).Value = "BCS"

# Row 5
$ws.Cells.Item(5,1).Value = "TUE"
$ws.Cells.Item(5,2).Value = "9:30-11:30"
$ws.Cells.Item(5,3).Value = "5CS024"
$ws.Cells.Item(5,4).Value = "Collaborative Development"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = "Lecture"
$ws.Cells.Item(5,7).Value = "Mr. Udaya Kandel"
$ws.Cells.Item(5,8).Value = "LT-03 Walsall"
$ws.Cells.Item(5,9).Value = "WLV"
$ws.Cells.Item(5,10).Value = "L5CG(9+10+11)"
$ws.Cells.Item(5,11).Value = 5
$ws.Cells.Item(5,12).Value = "BCS"

# Row 6
$ws.Cells.Item(6,1).Value = "TUE"
$ws.Cells.Item(6,2).Value = "12:00-14:00"
$ws.Cells.Item(6,3).Value = "5CS022"
$ws.Cells.Item(6,4).Value = "Human Computer Interaction"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = "Lecture"
$ws.Cells.Item(6,7).Value = "Mr. Pravash Karki"
$ws.Cells.Item(6,8).Value = "LT-02 Telford"
$ws.Cells.Item(6,9).Value = "WLV"
$ws.Cells.Item(6,10).Value = "L5CG(9+10+11)"
$ws.Cells.Item(6,11).Value = 5
$ws.Cells.Item(6,12).Value = "BCS"

# Row 7
$ws.Cells.Item(7,1).Value = "WED"
$ws.Cells.Item(7,2).Value = "9:00-11:00"
$ws.Cells.Item(7,3).Value = "5CS020"
$ws.Cells.Item(7,4).Value = "Distributed and Cloud Systems Programming"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = "Lecture"
$ws.Cells.Item(7,7).Value = "Mr. Sumanta Silwal"
$ws.Cells.Item(7,8).Value = "LT-02 Telford"
$ws.Cells.Item(7,9).Value = "WLV"
$ws.Cells.Item(7,10).Value = "L5CG(9+10+11)"
$ws.Cells.Item(7,11).Value = 5
$ws.Cells.Item(7,12).Value = "BCS"

# Row 8
$ws.Cells.Item(8,1).Value = "WED"
$ws.Cells.Item(8,2).Value = "12:30-14:30"
$ws.Cells.Item(8,3).Value = "5CS024"
$ws.Cells.Item(8,4).Value = "Collaborative Development"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = "Tutorial"
$ws.Cells.Item(8,7).Value = "Mr. Anmol Adhikari"
$ws.Cells.Item(8,8).Value = "TR-09  Chandragiri"
$ws.Cells.Item(8,9).Value = "HCK"
$ws.Cells.Item(8,10).Value = "L5CG10"
$ws.Cells.Item(8,11).Value = 5
$ws.Cells.Item(8,12).Value = "BCS"

# Row 9
$ws.Cells.Item(9,1).Value = "THU"
$ws.Cells.Item(9,2).Value = "13:00-15:00"
$ws.Cells.Item(9,3).Value = "5CS020"
$ws.Cells.Item(9,4).Value = "Distributed and Cloud Systems Programming"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = "Tutorial"
$ws.Cells.Item(9,7).Value = "Mr. Deepson Shrestha"
$ws.Cells.Item(9,8).Value = "SR-01 Bantok"
$ws.Cells.Item(9,9).Value = "WLV"
$ws.Cells.Item(9,10).Value = "L5CG10"
$ws.Cells.Item(9,11).Value = 5
$ws.Cells.Item(9,12).Value = "BCS"

# Row 10
$ws.Cells.Item(10,1).Value = "FRI"
$ws.Cells.Item(10,2).Value = "9:30-12:00"
$ws.Cells.Item(10,3).Value = "5CS020"
$ws.Cells.Item(10,4).Value = "Distributed and Cloud Systems Programming"
$ws.Cells.Item(10,5).Value = 2.5
$ws.Cells.Item(10,6).Value = "Workshop"
$ws.Cells.Item(10,7).Value = "Mr. Deepson Shrestha"
$ws.Cells.Item(10,8).Value = "Lab-02 Moseley"
$ws.Cells.Item(10,9).Value = "WLV"
$ws.Cells.Item(10,10).Value = "L5CG10"
$ws.Cells.Item(10,11).Value = 5
$ws.Cells.Item(10,12).Value = "BCS"
